$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A69").Value = "pi_obs"
$ws.Range("D69").Value = "ΔLN(GDPDEF)"

$ws.Range("A70").Value = "q_psi_obs"
$ws.Range("B70").Value = "SP500"
$ws.Range("D70").Value = "ΔLN(SP500/GDPDEF)"

$ws.Range("A71").Value = "l_h_obs"
$ws.Range("B71").Value = "Loans to households"
$ws.Range("D74").Value = "MORTG/100"

$ws.Range("A72").Value = "l_e_obs"
$ws.Range("B72").Value = "Loans to entrepreneurs"
$ws.Range("D72").Value = "ΔLN(TLBSNNCB/GDPDEF)"

$ws.Range("A73").Value = "d_obs"
$ws.Range("B73").Value = "Deposits"
$ws.Range("D73").Value = "ΔLN(DABSHNO/GDPDEF)"

$ws.Range("A74").Value = "i_h_obs"
$ws.Range("B74").Value = "Loan rate to households"
$ws.Range("D71").Value = "ΔLN(TLBSHNO/GDPDEF)"

$ws.Range("A75").Value = "i_e_obs"
$ws.Range("B75").Value = "Loan rate to entrepreneurs"
$ws.Range("D75").Value = "BAA/100"

$ws.Range("B69").Value = "Inflation"
